# Se agrego la implementacion de ColaEstatica y su clase Test.
# This script fills in the "Preparacion de la Prueba" timing row and the two
# new increment rows ("Crear interfaz Cola" / "Implementar Cola Estatica")
# on the "Metricas" sheet, matching the authored workbook edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metricas")

# --- Row 2: "Preparacion de la Prueba" timing -----------------------------
# B2 = tiempo real (duration), D2 = hora de inicio, E2 = hora de fin
# (C2 = E2-D2 is already a formula on the sheet and recalculates itself)
$ws.Range("B2").Value = [Math]::Round(5/1440, 15)          # 00:05:00 duration
$ws.Range("D2").Value = [Math]::Round(11.5/24, 15)         # 11:30:00
$ws.Range("E2").Value = [Math]::Round((11 + 35/60)/24, 15) # 11:35:00

# --- Row 6: "Crear interfaz Cola" increment -------------------------------
$ws.Range("A6").Value = "Crear interfaz Cola"
$ws.Range("B6").Value = 7
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = [Math]::Round(10/1440, 15)           # 00:10:00
$ws.Range("E6").Value = [Math]::Round((11 + 40/60)/24, 15)   # 11:40:00
$ws.Range("F6").Value = [Math]::Round((11 + 45/60)/24, 15)   # 11:45:00
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0

# --- Row 7: "Implementar Cola Estatica" increment -------------------------
$ws.Range("A7").Value = "Implementar Cola Estatica"
$ws.Range("B7").Value = 30
$ws.Range("C7").Value = 41
$ws.Range("D7").Value = [Math]::Round(20/1440, 15)           # 00:20:00
$ws.Range("E7").Value = [Math]::Round((11 + 47/60)/24, 15)   # 11:47:00
$ws.Range("F7").Value = [Math]::Round(12.5/24, 15)           # 12:30:00
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = [Math]::Round(30/1440, 15)           # 00:30:00

# Leave the final selection on B10, like the author's saved workbook.
$ws.Range("B10").Select()

$wb.Application.CalculateFull()
